$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.071262359619141
$ws.Range("B1").Value = 4.63306999206543
$ws.Range("C1").Value = 4.002201557159424
$ws.Range("D1").Value = 4.956319808959961
$ws.Range("E1").Value = 4.859043598175049
